$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.105.21"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -2.20%  '

$ws.Range('D3').Value = "'1.850.25"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.91%  '

$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').Value = "'0.6936"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.88%  '

$ws.Range('D6').Value = "'237.95"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.02%  '

$ws.Range('D7').Value = "'1.001"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').Value = "'0.07702"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +8.76%  '

$ws.Range('D9').Value = "'0.3032"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.77%  '

$ws.Range('D10').Value = "'23.24"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.22%  '

$ws.Range('D11').Value = "'0.08109"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.10%  '

$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = "'0.7257"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.86%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = "'1.830.26"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.42%  '

$ws.Range('D14').Value = "'5.209"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.80%  '

$ws.Range('D15').Value = "'88.99"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.21%  '

$ws.Range('D16').Value = "'29.103.80"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.29%  '

$ws.Range('D17').Value = "'5.749"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.13%  '

$ws.Range('D18').Value = "'0.000007775"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.19%  '

$ws.Range('D19').Value = "'13.19"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.06%  '

$ws.Range('D20').Value = "'236.39"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.45%  '

$ws.Range('D21').Value = "'1.000"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.04%  '

$ws.Range('D22').Value = "'2.094.60"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.72%  '

$ws.Range('D23').Value = "'1.001"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.19%  '

$ws.Range('D24').Value = "'7.601"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.60%  '

$ws.Range('E25').Value = '  -2.00%  '

$ws.Range('D26').Value = "'161.10"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.24%  '

$ws.Range('D27').Value = "'0.1432"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -6.67%  '

$ws.Range('D28').Value = "'18.04"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.30%  '

$ws.Range('E29').Value = '  -0.89%  '

$ws.Range('D30').Value = "'1.399"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.79%  '

$ws.Range('D31').Value = "'4.496"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.34%  '

$ws.Range('D32').Value = "'1.485"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.07%  '

$ws.Range('D33').Value = "'4.012"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.50%  '

$ws.Range('D34').Value = "'0.05224"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.71%  '

$ws.Range('D35').Value = "'1.184"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.65%  '

$ws.Range('E36').Value = '  +2.18%  '

$ws.Range('D37').Value = "'0.7001"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.93%  '

$ws.Range('D38').Value = "'2.658"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.11%  '

$ws.Range('D39').Value = "'0.01849"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.01%  '

$ws.Range('D40').Value = "'2.680"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.76%  '

$ws.Range('D41').Value = "'0.9195"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +6.30%  '

$ws.Range('D42').Value = "'6.005"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.37%  '

$ws.Range('D43').Value = "'1.080.75"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.88%  '

$ws.Range('E44').Value = '  -4.27%  '

$ws.Range('D45').Value = "'70.28"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.04%  '

$ws.Range('D46').Value = "'1.001"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.16%  '

$ws.Range('D47').Value = "'103.31"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.45%  '

$ws.Range('D48').Value = "'1.772"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.29%  '

$ws.Range('D49').Value = "'1.989.35"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.90%  '

$ws.Range('D50').Value = "'9.134"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.41%  '

$ws.Range('D51').Value = "'7.006"
$ws.Range('D51').ClearFormats()
